$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 277, pushing the existing data
# (old rows 277-281) down to become rows 279-283.
$ws.Rows.Item(277).Insert()
$ws.Rows.Item(277).Insert()

# New row 277: Early Majestic / Primera
$ws.Range("A277").Value = 7
$ws.Range("B277").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C277").Value = "Ñuble"
$ws.Range("D277").Value = 44890
$ws.Range("E277").Value = 16
$ws.Range("F277").Value = "Fruta"
$ws.Range("G277").Value = 100103
$ws.Range("H277").Value = "Frutos de hueso (carozo)"
$ws.Range("I277").Value = 100103004
$ws.Range("J277").Value = "Durazno"
$ws.Range("K277").Value = "Early Majestic"
$ws.Range("L277").Value = "Primera"
$ws.Range("M277").Value = 120
$ws.Range("N277").Value = 19000
$ws.Range("O277").Value = 20000
$ws.Range("P277").Value = 19500
$ws.Range("Q277").Value = "$/caja 16 kilos granel"
$ws.Range("R277").Value = "Región de O'Higgins"
$ws.Range("S277").Value = 1219
$ws.Range("T277").Value = 16

# New row 278: Early Majestic / Segunda
$ws.Range("A278").Value = 7
$ws.Range("B278").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C278").Value = "Ñuble"
$ws.Range("D278").Value = 44890
$ws.Range("E278").Value = 16
$ws.Range("F278").Value = "Fruta"
$ws.Range("G278").Value = 100103
$ws.Range("H278").Value = "Frutos de hueso (carozo)"
$ws.Range("I278").Value = 100103004
$ws.Range("J278").Value = "Durazno"
$ws.Range("K278").Value = "Early Majestic"
$ws.Range("L278").Value = "Segunda"
$ws.Range("M278").Value = 60
$ws.Range("N278").Value = 15000
$ws.Range("O278").Value = 15000
$ws.Range("P278").Value = 15000
$ws.Range("Q278").Value = "$/caja 16 kilos granel"
$ws.Range("R278").Value = "Región de O'Higgins"
$ws.Range("S278").Value = 938
$ws.Range("T278").Value = 16
